$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "WMT_Extract_SA"

$newSheet.Range("A1").Value = "Case_Ref_No"
$newSheet.Range("B1").Value = "Tier_Code"
$newSheet.Range("C1").Value = "Team_Code"
$newSheet.Range("D1").Value = "OM_Grade_Code"
$newSheet.Range("E1").Value = "OM_Key"
$newSheet.Range("F1").Value = "Location"
$newSheet.Range("G1").Value = "Disposal_Type_Desc"
$newSheet.Range("H1").Value = "Disposal_Type_Code"
$newSheet.Range("I1").Value = "Standalone_Order"

$colorVal = 0xCC + (0xCC * 256) + (0xFF * 65536)

$a1 = $newSheet.Range("A1")
$a1.Borders.LineStyle = 1
$a1.Borders.Color = $colorVal

foreach ($col in @("B","C","D","E","F","G","H","I")) {
    $cell = $newSheet.Range($col + "1")
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Color = $colorVal
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Color = $colorVal
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Color = $colorVal
}

Write-Output "done"
